# The presentation's main Design (slide master) was switched from the
# "Integral" (Red Violet) theme colours to the default Office Theme
# colour scheme. Font scheme and format scheme are unchanged -- only
# the theme colour scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink)
# differs between the two themes, so re-point each swatch individually
# via the Theme's ThemeColorScheme, exactly as PowerPoint's Design ->
# Colors picker (or VBA ThemeColorScheme.Colors(...).RGB) would.

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Colors([Microsoft.Office.Core.MsoThemeColorSchemeIndex]::msoThemeDark1).RGB = 0
$tcs.Colors([Microsoft.Office.Core.MsoThemeColorSchemeIndex]::msoThemeLight1).RGB = 16777215
$tcs.Colors([Microsoft.Office.Core.MsoThemeColorSchemeIndex]::msoThemeDark2).RGB = 6968388
$tcs.Colors([Microsoft.Office.Core.MsoThemeColorSchemeIndex]::msoThemeLight2).RGB = 15132391
$tcs.Colors([Microsoft.Office.Core.MsoThemeColorSchemeIndex]::msoThemeAccent1).RGB = 13998939
$tcs.Colors([Microsoft.Office.Core.MsoThemeColorSchemeIndex]::msoThemeAccent2).RGB = 3243501
$tcs.Colors([Microsoft.Office.Core.MsoThemeColorSchemeIndex]::msoThemeAccent3).RGB = 10855845
$tcs.Colors([Microsoft.Office.Core.MsoThemeColorSchemeIndex]::msoThemeAccent4).RGB = 49407
$tcs.Colors([Microsoft.Office.Core.MsoThemeColorSchemeIndex]::msoThemeAccent5).RGB = 12874308
$tcs.Colors([Microsoft.Office.Core.MsoThemeColorSchemeIndex]::msoThemeAccent6).RGB = 4697456
$tcs.Colors([Microsoft.Office.Core.MsoThemeColorSchemeIndex]::msoThemeHyperlink).RGB = 12673797
$tcs.Colors([Microsoft.Office.Core.MsoThemeColorSchemeIndex]::msoThemeFollowedHyperlink).RGB = 7491477
